{"js": "// Office.js (Word JavaScript API) script implementing the edit described by the diff:\n//  1. Remove the \"License Information\" (Heading2) paragraph.\n//  2. Rewrite the license paragraph (\"... is based on: unfoldingWord\u00ae Translation\n//     Words, unfoldingWord, 2022, which is licensed under a CC BY-SA 4.0 license.\")\n//     as new copyright/adaptation text, dropping both hyperlinks (now plain text)\n//     and merging away the following \"This PDF version is provided under the\n//     same license.\" paragraph entirely.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// --- Step 1: locate + delete the \"License Information\" heading paragraph ---\nlet licenseInfoPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"License Information\") {\n    licenseInfoPara = paragraphs.items[i];\n    break;\n  }\n}\nif (licenseInfoPara) {\n  licenseInfoPara.delete();\n  await context.sync();\n}\n\n// --- Step 2: locate the license paragraph + the paragraph right after it ---\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\n\nlet licenseParaIndex = -1;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text.indexOf(\"is based on\") >= 0) {\n    licenseParaIndex = i;\n    break;\n  }\n}\n\nif (licenseParaIndex >= 0) {\n  const licensePara = paragraphs2.items[licenseParaIndex];\n  const nextPara = paragraphs2.items[licenseParaIndex + 1];\n\n  // Range spanning both paragraphs (the license paragraph + the \"This PDF\n  // version...\" paragraph that follows it), so replacing it merges the two\n  // into a single new paragraph (i.e. the second paragraph disappears).\n  let targetRange;\n  if (nextPara) {\n    targetRange = licensePara.getRange(\"Start\").expandTo(nextPara.getRange(\"End\"));\n  } else {\n    targetRange = licensePara.getRange();\n  }\n\n  const newParaOoxml = `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">` +\n    `<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">` +\n    `<pkg:xmlData>` +\n    `<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">` +\n    `<w:body>` +\n    `<w:p>` +\n      `<w:pPr/>` +\n      `<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr></w:r>` +\n      `<w:r><w:rPr><w:b/><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr><w:t>unfoldingWord\\u00AE Translation Words</w:t></w:r>` +\n      `<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr><w:t xml:space=\"preserve\"> \\u00A9 2022 unfoldingWord. Released under CC BY-SA 4.0 license. </w:t></w:r>` +\n      `<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr><w:t>unfoldingWord\\u00AE Translation Words</w:t></w:r>` +\n      `<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr><w:t xml:space=\"preserve\"> has been adapted in the following languages: Tok Pisin, Arabic (\\u0639\\u0631\\u0628\\u064A), French (Fran\\u00E7ais), Hindi (\\u0939\\u093F\\u0902\\u0926\\u0940), Indonesian (Bahasa Indonesia), Portuguese (Portugu\\u00EAs), Russian (\\u0420\\u0443\\u0441\\u0441\\u043A\\u0438\\u0439), Spanish (Espa\\u00F1ol), Swahili (Kiswahili), and Simplified Chinese (\\u7B80\\u4F53\\u4E2D\\u6587) from </w:t></w:r>` +\n      `<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr><w:t>unfoldingWord\\u00AE Translation Words</w:t></w:r>` +\n      `<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr><w:t xml:space=\"preserve\"> \\u00A9 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual</w:t></w:r>` +\n      `<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr></w:r>` +\n      `<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr></w:r>` +\n    `</w:p>` +\n    `</w:body>` +\n    `</w:document>` +\n    `</pkg:xmlData>` +\n    `</pkg:part>` +\n    `</pkg:package>`;\n\n  targetRange.insertOoxml(newParaOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script implementing the edit described\n# by the diff:\n#   1. Remove the \"License Information\" (Heading2) paragraph.\n#   2. Rewrite the license paragraph (\"... is based on: unfoldingWord(R)\n#      Translation Words, unfoldingWord, 2022, which is licensed under a\n#      CC BY-SA 4.0 license.\") as new copyright/adaptation text, dropping\n#      both hyperlinks (now plain text) and merging away the following\n#      \"This PDF version is provided under the same license.\" paragraph\n#      entirely.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: locate + delete the \"License Information\" heading paragraph ---\n$licenseInfoPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"License Information\") {\n        $licenseInfoPara = $p\n        break\n    }\n}\nif ($licenseInfoPara -ne $null) {\n    $licenseInfoPara.Range.Delete()\n}\n\n# --- Step 2: locate the license paragraph + the paragraph right after it ---\n$licenseParaIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t.IndexOf(\"is based on\") -ge 0) {\n        $licenseParaIndex = $i\n        break\n    }\n}\n\nif ($licenseParaIndex -ge 0) {\n    $licensePara = $d.Paragraphs.Item($licenseParaIndex)\n    $targetRange = $null\n    if ($licenseParaIndex -lt $d.Paragraphs.Count) {\n        $nextPara = $d.Paragraphs.Item($licenseParaIndex + 1)\n        $targetRange = $d.Range($licensePara.Range.Start, $nextPara.Range.End)\n    } else {\n        $targetRange = $licensePara.Range\n    }\n\n    $newParaOoxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n        '<w:p>' +\n          '<w:pPr/>' +\n          '<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr></w:r>' +\n          '<w:r><w:rPr><w:b/><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr><w:t>unfoldingWord&#174; Translation Words</w:t></w:r>' +\n          '<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr><w:t xml:space=\"preserve\"> &#169; 2022 unfoldingWord. Released under CC BY-SA 4.0 license. </w:t></w:r>' +\n          '<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr><w:t>unfoldingWord&#174; Translation Words</w:t></w:r>' +\n          '<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr><w:t xml:space=\"preserve\"> has been adapted in the following languages: Tok Pisin, Arabic (&#1593;&#1585;&#1576;&#1610;), French (Fran&#231;ais), Hindi (&#2361;&#2367;&#2306;&#2342;&#2368;), Indonesian (Bahasa Indonesia), Portuguese (Portugu&#234;s), Russian (&#1056;&#1091;&#1089;&#1089;&#1082;&#1080;&#1081;), Spanish (Espa&#241;ol), Swahili (Kiswahili), and Simplified Chinese (&#31616;&#20307;&#20013;&#25991;) from </w:t></w:r>' +\n          '<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr><w:t>unfoldingWord&#174; Translation Words</w:t></w:r>' +\n          '<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr><w:t xml:space=\"preserve\"> &#169; 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual</w:t></w:r>' +\n          '<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr></w:r>' +\n          '<w:r><w:rPr><w:lang w:val=\"hi_IN\" w:bidi=\"hi_IN\"/></w:rPr></w:r>' +\n        '</w:p>' +\n        '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n\n    $targetRange.InsertXML($newParaOoxml)\n}\n"}
